$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3907.5293
$ws.Range("I100").Value = 2080.6
$ws.Range("K100").Value = 2080.6
$ws.Range("M100").Value = -1539.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1885488.2
$ws.Range("J112").Value = 1993161.9
$ws.Range("L112").Value = 5979485.699999999
$ws.Range("N112").Value = -5981701.699999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4983
$ws.Range("I116").Value = 4815.8184
$ws.Range("K116").Value = 4815.8184
$ws.Range("M116").Value = -1373.8184

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5949480
$ws.Range("I137").Value = 9192217
$ws.Range("K137").Value = 27576651
$ws.Range("M137").Value = -27574101

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1665.3334
$ws.Range("I102").Value = 1614.8182
$ws.Range("K102").Value = 1614.8182
$ws.Range("M102").Value = 7.181800000000067

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5493.5713
$ws.Range("I20").Value = 5409.1665
$ws.Range("K20").Value = 5409.1665
$ws.Range("M20").Value = -5162.1665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5273.5
$ws.Range("I99").Value = 7500
$ws.Range("K99").Value = 7500
$ws.Range("M99").Value = -6002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2458.111
$ws.Range("I105").Value = 2308.5881
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 2308.5881
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -561.5880999999999
$ws.Range("N105").Value = -8494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 712.5
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3654
$ws.Range("I58").Value = 1617.4
$ws.Range("J58").Value = 6199.75
$ws.Range("K58").Value = 1617.4
$ws.Range("L58").Value = 6199.75
$ws.Range("M58").Value = -1414.4
$ws.Range("N58").Value = -6605.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 28995
$ws.Range("J106").Value = 28995
$ws.Range("L106").Value = 28995
$ws.Range("N106").Value = -31519

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3654
$ws.Range("I136").Value = 1617.4
$ws.Range("J136").Value = 6199.75
$ws.Range("K136").Value = 4852.200000000001
$ws.Range("L136").Value = 18599.25
$ws.Range("M136").Value = -2302.200000000001
$ws.Range("N136").Value = -23699.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 201247
$ws.Range("I109").Value = 201247
$ws.Range("K109").Value = 603741
$ws.Range("M109").Value = -602701

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 5624.3667
$ws.Range("I132").Value = 7130.6665
$ws.Range("J132").Value = 2109.6667
$ws.Range("K132").Value = 64175.9985
$ws.Range("L132").Value = 18987.0003
$ws.Range("M132").Value = -61645.9985
$ws.Range("N132").Value = -24047.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1770.4642
$ws.Range("I97").Value = 1654.2778
$ws.Range("J97").Value = 1979.6
$ws.Range("K97").Value = 1654.2778
$ws.Range("L97").Value = 1979.6
$ws.Range("M97").Value = -1158.2778
$ws.Range("N97").Value = -2971.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1386.1818
$ws.Range("I107").Value = 875.2857
$ws.Range("J107").Value = 2280.25
$ws.Range("K107").Value = 875.2857
$ws.Range("L107").Value = 2280.25
$ws.Range("M107").Value = 1044.7143
$ws.Range("N107").Value = -6120.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 15353.637
$ws.Range("I113").Value = 21213
$ws.Range("K113").Value = 21213
$ws.Range("M113").Value = -19043

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3476.889
$ws.Range("I122").Value = 3465.5
$ws.Range("J122").Value = 3499.6667
$ws.Range("K122").Value = 10396.5
$ws.Range("L122").Value = 10499.0001
$ws.Range("M122").Value = -7946.5
$ws.Range("N122").Value = -15399.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3429.5
$ws.Range("I46").Value = 3020.4211
$ws.Range("J46").Value = 4984
$ws.Range("K46").Value = 3020.4211
$ws.Range("L46").Value = 4984
$ws.Range("M46").Value = -2832.4211
$ws.Range("N46").Value = -5360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4985.4287
$ws.Range("I93").Value = 1724.5
$ws.Range("J93").Value = 9333.333000000001
$ws.Range("K93").Value = 1724.5
$ws.Range("L93").Value = 9333.333000000001
$ws.Range("M93").Value = -476.5
$ws.Range("N93").Value = -11829.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 9254.6
$ws.Range("I100").Value = 8987.395500000001
$ws.Range("K100").Value = 8987.395500000001
$ws.Range("M100").Value = -8446.395500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 40000
$ws.Range("J47").Value = 40000
$ws.Range("L47").Value = 40000
$ws.Range("N47").Value = -41144

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 49671.332
$ws.Range("J48").Value = 49671.332
$ws.Range("L48").Value = 49671.332
$ws.Range("N48").Value = -50809.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1346.6364
$ws.Range("I100").Value = 1281.6
$ws.Range("K100").Value = 2563.2
$ws.Range("M100").Value = -2022.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 45806.75
$ws.Range("J105").Value = 45806.75
$ws.Range("L105").Value = 45806.75
$ws.Range("N105").Value = -52794.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3219.125
$ws.Range("I122").Value = 3219.125
$ws.Range("K122").Value = 9657.375
$ws.Range("M122").Value = -7207.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 59950
$ws.Range("J123").Value = 59950
$ws.Range("L123").Value = 59950
$ws.Range("N123").Value = -69750

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3134.2144
$ws.Range("I126").Value = 3117.5
$ws.Range("J126").Value = 3156.5
$ws.Range("K126").Value = 9352.5
$ws.Range("L126").Value = 9469.5
$ws.Range("M126").Value = -6882.5
$ws.Range("N126").Value = -14409.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1617.5111
$ws.Range("I132").Value = 1576.4651
$ws.Range("K132").Value = 4729.3953
$ws.Range("M132").Value = -2199.3953
